# Append 9 new data rows (58-66) to Sheet1, matching the source diff that
# extends the used range from A1:I57 to A1:I66.
#
# Columns: A=id, B=patientId, C=name, D=age, E=audioResult, F=videoResult,
#          G=timestamp, H=mmseScore, I=status

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Reuse the existing date/time number format already applied to the
# "timestamp" column (G) so the new cells share the same style index
# instead of Excel minting a brand-new one.
$tsFormat = $ws.Cells.Item(57, 7).NumberFormat()

$rows = @(
    @{ Row = 58; Id = 98;  PatientId = "PT301224163607"; Name = "dfdgg"; Age = 34; Audio = $null; Video = $null;    Timestamp = 45656.69174768519; Mmse = 0; Status = "Severe" },
    @{ Row = 59; Id = 99;  PatientId = "PT301224163630"; Name = "gg";    Age = 55; Audio = $null; Video = $null;    Timestamp = 45656.69201388889; Mmse = 0; Status = "Severe" },
    @{ Row = 60; Id = 100; PatientId = "PT301224163705"; Name = "qq";    Age = 1;  Audio = $null; Video = $null;    Timestamp = 45656.69241898148; Mmse = 0; Status = "Severe" },
    @{ Row = 61; Id = 101; PatientId = "PT301224163741"; Name = "ww";    Age = 2;  Audio = $null; Video = $null;    Timestamp = 45656.69283564815; Mmse = 0; Status = "Severe" },
    @{ Row = 62; Id = 102; PatientId = "PT301224164034"; Name = "q";     Age = 1;  Audio = $null; Video = $null;    Timestamp = 45656.69483796296; Mmse = 0; Status = "Severe" },
    @{ Row = 63; Id = 103; PatientId = "PT301224165257"; Name = "qa";    Age = 23; Audio = $null; Video = $null;    Timestamp = 45656.7034375;      Mmse = 0; Status = "Severe" },
    @{ Row = 64; Id = 104; PatientId = "PT301224165412"; Name = "qq";    Age = 22; Audio = $null; Video = $null;    Timestamp = 45656.70430555556; Mmse = 0; Status = "Severe" },
    @{ Row = 65; Id = 105; PatientId = "PT301224165530"; Name = "az";    Age = 45; Audio = "mci"; Video = "Normal"; Timestamp = 45656.70572916666; Mmse = 0; Status = "Severe" },
    @{ Row = 66; Id = 108; PatientId = "PT301224170458"; Name = "rt";    Age = 6;  Audio = "mci"; Video = "Dementia"; Timestamp = 45656.71229166666; Mmse = 0; Status = "Severe" }
)

foreach ($r in $rows) {
    $row = $r.Row
    $ws.Cells.Item($row, 1).Value = $r.Id
    $ws.Cells.Item($row, 2).Value = $r.PatientId
    $ws.Cells.Item($row, 3).Value = $r.Name
    $ws.Cells.Item($row, 4).Value = $r.Age
    if ($r.Audio -ne $null) {
        $ws.Cells.Item($row, 5).Value = $r.Audio
    }
    if ($r.Video -ne $null) {
        $ws.Cells.Item($row, 6).Value = $r.Video
    }
    $ws.Cells.Item($row, 7).Value = $r.Timestamp
    $ws.Cells.Item($row, 7).NumberFormat = $tsFormat
    $ws.Cells.Item($row, 8).Value = $r.Mmse
    $ws.Cells.Item($row, 9).Value = $r.Status
}
